# Update the DCL test case strings in Sheet1 to reflect the renamed DB
# column ("user" -> "name"), per commit message:
#   "change dcl case as the db column name changed by developer"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = "select host,name,AUTHENTICATION_STRING from mysql.user where name='dcltest1'"
$ws.Range("M3").Value = "select host,name,AUTHENTICATION_STRING from mysql.user where name='dcltest2'"
$ws.Range("M4").Value = "select host,name,AUTHENTICATION_STRING from mysql.user where name='dcltest3'"
$ws.Range("M5").Value = "select name,AUTHENTICATION_STRING from mysql.user where name='dcltest4'"

# L5 text is unchanged, but keep it explicit/consistent in case it was touched.
$ws.Range("L5").Value = "create user 'dcltest4'@'localhost' identified by 'abc123'"

# Update the selected/active cell shown in the sheet view.
$ws.Range("M2").Select()
